# Daily attendance processing
# Normalizes the "Recorded By" (column G) values so that the automated
# "System"/"system" recorder entry is always listed first, followed by the
# human recorders/emails, preserving their original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $orig = $cell.Value2

    if ($orig -eq $null) { continue }
    if ($orig -eq "") { continue }

    $parts = $orig.Split(",")

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        $pt = $p.Trim()
        if ($pt.ToLower() -eq "system") {
            $systemParts += $pt
        } else {
            $otherParts += $pt
        }
    }

    $newParts = $systemParts + $otherParts
    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -ne $orig) {
        $cell.Value2 = $newVal
    }
}
